$d = $word.ActiveDocument

# 1) SEXO_5 -> SEXO_14 in the "MEXICAN{{SEXO_5}}" placeholder
$d.Content.Find.Execute(
    "MEXICAN{{SEXO_5}}", $true, $false, $false, $false, $false,
    $true, 1, $false, "MEXICAN{{SEXO_14}}", 2)

# 2) " ... SOBRE LA ESCRITURACIÓN POR CUESTIONES AJENAS A ÉSTE E." ->
#    " ... SOBRE LA COMPRAVENTA POR CUESTIONES AJENAS A ÉSTE."
$d.Content.Find.Execute(
    "ESCRITURACIÓN POR CUESTIONES AJENAS A ÉSTE E.", $true, $false, $false, $false, $false,
    $true, 1, $false, "COMPRAVENTA POR CUESTIONES AJENAS A ÉSTE.", 2)
